$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.565.55"
$ws.Range("E2").Value = "'  +2.55%  "
$ws.Range("D3").Value = "'1.669.83"
$ws.Range("E3").Value = "'  +2.04%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("D5").Value = "'239.22"
$ws.Range("E5").Value = "'  +1.56%  "
$ws.Range("E6").Value = "'  -0.02%  "
$ws.Range("D7").Value = "'0.4783"
$ws.Range("E7").Value = "'  +1.07%  "
$ws.Range("D8").Value = "'0.2628"
$ws.Range("E8").Value = "'  +2.83%  "
$ws.Range("D9").Value = "'0.06171"
$ws.Range("E9").Value = "'  +2.83%  "
$ws.Range("D10").Value = "'1.667.39"
$ws.Range("E10").Value = "'  +1.80%  "
$ws.Range("D11").Value = "'0.07001"
$ws.Range("E11").Value = "'  -2.31%  "
$ws.Range("D12").Value = "'14.88"
$ws.Range("E12").Value = "'  +1.07%  "
$ws.Range("D13").Value = "'0.5897"
$ws.Range("E13").Value = "'  -4.06%  "
$ws.Range("D14").Value = "'4.380"
$ws.Range("E14").Value = "'  -0.51%  "
$ws.Range("D15").Value = "'75.28"
$ws.Range("E15").Value = "'  +3.84%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "'  -0.04%  "
$ws.Range("D17").Value = "'1.0000"
$ws.Range("E17").Value = "'  +0.25%  "
$ws.Range("D18").Value = "'25.549.90"
$ws.Range("E18").Value = "'  +2.55%  "
$ws.Range("D19").Value = "'0.000006751"
$ws.Range("E19").Value = "'  +2.93%  "
$ws.Range("D20").Value = "'11.44"
$ws.Range("E20").Value = "'  +2.08%  "
$ws.Range("D21").Value = "'1.882.83"
$ws.Range("E21").Value = "'  +2.01%  "
$ws.Range("D22").Value = "'4.437"
$ws.Range("E22").Value = "'  +0.72%  "
$ws.Range("D23").Value = "'8.753"
$ws.Range("E23").Value = "'  +2.22%  "
$ws.Range("D24").Value = "'5.280"
$ws.Range("E24").Value = "'  +0.44%  "
$ws.Range("D25").Value = "'136.69"
$ws.Range("E25").Value = "'  +3.22%  "
$ws.Range("D26").Value = "'15.04"
$ws.Range("E26").Value = "'  +1.85%  "
$ws.Range("E27").Value = "'  +1.59%  "
$ws.Range("D28").Value = "'1.722"
$ws.Range("E28").Value = "'  +4.22%  "
$ws.Range("D29").Value = "'104.78"
$ws.Range("E29").Value = "'  +2.31%  "
$ws.Range("D30").Value = "'3.963"
$ws.Range("E30").Value = "'  +6.61%  "
$ws.Range("D31").Value = "'0.07829"
$ws.Range("E31").Value = "'  +1.07%  "
$ws.Range("D32").Value = "'3.650"
$ws.Range("E32").Value = "'  +3.37%  "
$ws.Range("D33").Value = "'0.9991"
$ws.Range("E33").Value = "'  -0.06%  "
$ws.Range("D34").Value = "'0.04227"
$ws.Range("E34").Value = "'  -3.34%  "
$ws.Range("D35").Value = "'2.622"
$ws.Range("E35").Value = "'  +0.90%  "
$ws.Range("B36").Value = "'ARBITRUM"
$ws.Range("C36").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9566"
$ws.Range("E36").Value = "'  +4.48%  "
$ws.Range("B37").Value = "'ImmutableX"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.6097"
$ws.Range("E37").Value = "'  +5.17%  "
$ws.Range("D38").Value = "'2.594"
$ws.Range("E38").Value = "'  +2.41%  "
$ws.Range("D39").Value = "'0.8604"
$ws.Range("E39").Value = "'  +4.48%  "
$ws.Range("D40").Value = "'0.9999"
$ws.Range("E40").Value = "'  +0.15%  "
$ws.Range("D41").Value = "'1.871"
$ws.Range("E41").Value = "'  +4.41%  "
$ws.Range("D42").Value = "'0.01477"
$ws.Range("E42").Value = "'  -4.87%  "
$ws.Range("D43").Value = "'96.58"
$ws.Range("E43").Value = "'  -0.74%  "
$ws.Range("D44").Value = "'0.3769"
$ws.Range("E44").Value = "'  +1.73%  "
$ws.Range("D45").Value = "'4.873"
$ws.Range("E45").Value = "'  +2.86%  "
$ws.Range("D46").Value = "'0.1117"
$ws.Range("E46").Value = "'  -1.43%  "
$ws.Range("D47").Value = "'6.227"
$ws.Range("E48").Value = "'  +1.29%  "
$ws.Range("D49").Value = "'29.92"
$ws.Range("E49").Value = "'  +1.47%  "
$ws.Range("D50").Value = "'7.375"
$ws.Range("E50").Value = "'  +3.03%  "
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "'  +0.15%  "

Write-Output "Applied cryptos update"
